# Apply "Add: mas opciones de RET PER" changes to the Control sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (blue fill / white font / border) from an existing
# header cell so the newly created cells (J1..O1) pick it up too.
$ws.Range("A1").Copy()
$ws.Range("J1:O1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the existing SIRCREB / SIFERE headers two columns to the right
# (they now live at K1 / L1) and add the new ARBA / CABA-AGIP split headers.
$ws.Range("K1").Value = "SIRCREB"
$ws.Range("L1").Value = "SIFERE"
$ws.Range("M1").Value = "ARBA"
$ws.Range("N1").Value = "CABA - AGIP RET"
$ws.Range("O1").Value = "CABA - AGIP PER"

# Replace the old G1:I1 headers with the new PER headers, and add the new
# SUSS header right after them at J1.
$ws.Range("G1").Value = "PER 216"
$ws.Range("H1").Value = "PER 217"
$ws.Range("I1").Value = "PER 767"
$ws.Range("J1").Value = "SUSS"

# Keep the active selection in line with the authored workbook.
$ws.Range("J2").Select()
